$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells to text format so numeric-looking strings
# (e.g. "1.00", "56.537.97") are stored as text, matching the source data.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '56.537.97'
$ws.Range("E2").Value = '  +4.23%  '
$ws.Range("D3").Value = '3.000.26'
$ws.Range("E3").Value = '  +4.87%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '507.54'
$ws.Range("E5").Value = '  +9.00%  '
$ws.Range("D6").Value = '136.91'
$ws.Range("E6").Value = '  +11.00%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '0.432'
$ws.Range("E8").Value = '  +8.16%  '
$ws.Range("D9").Value = '7.58'
$ws.Range("E9").Value = '  +15.73%  '
$ws.Range("D10").Value = '0.108'
$ws.Range("E10").Value = '  +14.66%  '
$ws.Range("D11").Value = '0.351'
$ws.Range("E11").Value = '  +9.07%  '
$ws.Range("D12").Value = '0.129'
$ws.Range("E12").Value = '  +5.99%  '
$ws.Range("D13").Value = '3.519.42'
$ws.Range("E13").Value = '  +4.93%  '
$ws.Range("D14").Value = '25.42'
$ws.Range("E14").Value = '  +11.23%  '
$ws.Range("D15").Value = '0.0000154'
$ws.Range("E15").Value = '  +18.01%  '
$ws.Range("D16").Value = '56.631.75'
$ws.Range("E16").Value = '  +4.41%  '
$ws.Range("D17").Value = '3.004.88'
$ws.Range("E17").Value = '  +4.77%  '
$ws.Range("D18").Value = '5.79'
$ws.Range("E18").Value = '  +10.39%  '
$ws.Range("D19").Value = '12.42'
$ws.Range("E19").Value = '  +11.22%  '
$ws.Range("D20").Value = '7.81'
$ws.Range("E20").Value = '  +12.46%  '
$ws.Range("D21").Value = '326.46'
$ws.Range("E21").Value = '  +12.17%  '
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").Value = '0.478'
$ws.Range("E23").Value = '  +10.64%  '
$ws.Range("D24").Value = '62.43'
$ws.Range("E24").Value = '  +7.91%  '
$ws.Range("D25").Value = '0.169'
$ws.Range("E25").Value = '  +14.50%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("D27").Value = '0.0₃0912'
$ws.Range("E27").Value = '  +16.78%  '
$ws.Range("D28").Value = '6.52'
$ws.Range("E28").Value = '  +9.30%  '
$ws.Range("D29").Value = '7.01'
$ws.Range("E29").Value = '  +16.46%  '
$ws.Range("D30").Value = '1.25'
$ws.Range("E30").Value = '  +16.11%  '
$ws.Range("D31").Value = '1.77'
$ws.Range("E31").Value = '  +11.71%  '
$ws.Range("D32").Value = '20.63'
$ws.Range("E32").Value = '  +13.05%  '
$ws.Range("D33").Value = '155.68'
$ws.Range("E33").Value = '  +14.00%  '
$ws.Range("D34").Value = '4.49'
$ws.Range("E34").Value = '  +10.50%  '
$ws.Range("D35").Value = '5.60'
$ws.Range("E35").Value = '  +6.01%  '
$ws.Range("D36").Value = '1.26'
$ws.Range("E36").Value = '  +6.11%  '
$ws.Range("D37").Value = '0.0674'
$ws.Range("E37").Value = '  +11.25%  '
$ws.Range("D38").Value = '24.12'
$ws.Range("E38").Value = '  +6.61%  '
$ws.Range("D39").Value = '3.039.06'
$ws.Range("E39").Value = '  +5.21%  '
$ws.Range("D40").Value = '36.67'
$ws.Range("E40").Value = '  +4.64%  '
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("D42").Value = '0.650'
$ws.Range("E42").Value = '  +8.59%  '
$ws.Range("D43").Value = '2.267.86'
$ws.Range("E43").Value = '  +12.25%  '
$ws.Range("D44").Value = '0.996'
$ws.Range("E44").Value = '  +8.49%  '
$ws.Range("D45").Value = '1.40'
$ws.Range("E45").Value = '  +9.59%  '
$ws.Range("D46").Value = '3.60'
$ws.Range("E46").Value = '  +8.98%  '
$ws.Range("D47").Value = '1.97'
$ws.Range("E47").Value = '  +26.44%  '
$ws.Range("D48").Value = '0.0236'
$ws.Range("E48").Value = '  +12.18%  '
$ws.Range("D49").Value = '5.75'
$ws.Range("E49").Value = '  +9.36%  '
$ws.Range("D50").Value = '19.00'
$ws.Range("E50").Value = '  +9.12%  '
$ws.Range("D51").Value = '0.0870'
$ws.Range("E51").Value = '  +12.50%  '

# Restore default cell style on column D cells so no stray style index is left
# applied to these cells (keeps them visually/structurally identical to before).
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
